# Updated cryptos list on Fri Mar  1 10:24:58 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for the crypto rows.
# A leading apostrophe forces Excel to keep values that look like plain
# numbers (e.g. "407.92") stored as text, matching the original inline
# string cell type instead of letting COM auto-coerce them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.088.80"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.422.62"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'407.92"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'134.51"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("D7").Value = "'0.592"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "'0.121"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").Value = "'42.77"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'8.45"
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "'19.86"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "3.447.21"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "62.023.19"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "'11.04"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "'0.0000130"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'3.20"
$ws.Range("E20").Value = "  -5.27%  "
$ws.Range("D21").Value = "'84.50"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("D22").Value = "'314.00"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'12.85"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'4.79"
$ws.Range("E25").Value = "  +9.77%  "
$ws.Range("D26").Value = "'29.66"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "'7.65"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").Value = "'2.76"
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("D30").Value = "'0.174"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("D32").Value = "'42.99"
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  -6.53%  "
$ws.Range("D35").Value = "'0.0484"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").Value = "'51.89"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "'3.42"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'137.39"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "'0.298"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").Value = "'4.03"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  -5.96%  "
$ws.Range("D46").Value = "'2.23"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "'21.42"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").Value = "2.123.38"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("E49").Value = "  -5.24%  "
$ws.Range("D50").Value = "'1.91"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "'0.0351"
$ws.Range("E51").Value = "  +4.02%  "
